$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 18 new blank rows at row 507, pushing existing rows 507-521 down to 525-539
$ws.Range("A507:A524").EntireRow.Insert()

# Populate the new rows (key in column A, pt translation in column B) for svat rules 30-35
$ws.Range("A507").Value = 'svat_t30'
$ws.Range("B507").Value = 'Teste à utilização da conta 261- "Acionistas - Acionistas com subscrição"'
$ws.Rows.Item(507).RowHeight = 32

$ws.Range("A508").Value = 'svat_t30_ok'
$ws.Range("B508").Value = 'A conta 261 deve ser utilizada exclusivamente por sociedades anónimas e comandita por ações. Sem exceções.'
$ws.Rows.Item(508).RowHeight = 48

$ws.Range("A509").Value = 'svat_t30_nok'
$ws.Range("B509").Value = 'A conta 261 deve ser utilizada exclusivamente por sociedades anónimas e comandita por ações. Verifique o ''tipo de sociedade'' na ficha da empresa ou os movimentos nas seguintes contas:'
$ws.Rows.Item(509).RowHeight = 64

$ws.Range("A510").Value = 'svat_t31'
$ws.Range("B510").Value = 'Teste à utilização da conta 262- "Sócios - Quotas não liberadas"'
$ws.Rows.Item(510).RowHeight = 32

$ws.Range("A511").Value = 'svat_t31_ok'
$ws.Range("B511").Value = 'A conta 262 deve ser utilizada exclusivamente por sociedades anónimas e comandita por ações. Sem exceções.'
$ws.Rows.Item(511).RowHeight = 48

$ws.Range("A512").Value = 'svat_t31_nok'
$ws.Range("B512").Value = 'A conta 262 deve ser utilizada exclusivamente por sociedades anónimas e comandita por ações. Verifique o ''tipo de sociedade'' na ficha da empresa ou os movimentos nas seguintes contas:'
$ws.Rows.Item(512).RowHeight = 64

$ws.Range("A513").Value = 'svat_t32'
$ws.Range("B513").Value = 'Teste à utilização da conta 263-"Adiantamentos por conta de lucros"'
$ws.Rows.Item(513).RowHeight = 32

$ws.Range("A514").Value = 'svat_t32_ok'
$ws.Range("B514").Value = 'A conta 263 não deve ser utilizada por sociedades anónimas. Sem exceções.'
$ws.Rows.Item(514).RowHeight = 32

$ws.Range("A515").Value = 'svat_t32_nok'
$ws.Range("B515").Value = 'A conta 263 não deve ser utilizada por sociedades anónimas. Nas sociedades anónimas, os adiantamentos devem estar representados no saldo devedor da conta 89 - "Dividendos antecipados". Verifique o saldo das contas:'
$ws.Rows.Item(515).RowHeight = 80

$ws.Range("A516").Value = 'account_269'
$ws.Range("B516").Value = 'Teste à conta 269-  "Acionistas/sócios - Perdas por imparidade acumuladas"'
$ws.Rows.Item(516).RowHeight = 32

$ws.Range("A517").Value = 'account_269_ok'
$ws.Range("B517").Value = 'Foi efetuado um teste sobre a conta 269 - "Acionistas/sócios - Perdas por imparidade acumuladas"e verificado que o saldo credor desta conta é igual ou  inferior à soma algébrica dos saldos das contas 261;262;263;266;267 e 268. Teste realizado com sucesso.'
$ws.Rows.Item(517).RowHeight = 96

$ws.Range("A518").Value = 'account_269_nok'
$ws.Range("B518").Value = 'Foi efetuado um teste sobre a conta 269 - "Acionistas/sócios - Perdas por imparidade acumuladas"e verificado que o saldo credor desta conta é igual ou  superior à soma algébrica dos saldos das contas 261;262;263;266;267 e 268. Esta situação deve ser regularizada para uma correta submissão do ficheiro SAF-T.'
$ws.Rows.Item(518).RowHeight = 112

$ws.Range("A519").Value = 'account_269_var'
$ws.Range("B519").Value = 'Teste à conta 269-  "Acionistas/sócios - Perdas por imparidade acumuladas"'
$ws.Rows.Item(519).RowHeight = 32

$ws.Range("A520").Value = 'account_269_var_ok'
$ws.Range("B520").Value = 'Foi efetuado um teste sobre a conta 269 - "Acionistas/sócios - Perdas por imparidade acumuladas"e verificado que o saldo credor desta conta, entidade a entidade, é igual ou  inferior à soma algébrica dos saldos das contas 261;262;263;266;267 e 268. Teste realizado com sucesso.'
$ws.Rows.Item(520).RowHeight = 112

$ws.Range("A521").Value = 'account_269_var_nok'
$ws.Range("B521").Value = 'Foi efetuado um teste sobre a conta 269 - "Acionistas/sócios - Perdas por imparidade acumuladas"e verificado que o saldo credor desta conta, entidade a entidade, é igual ou  superior à soma algébrica dos saldos das contas 261;262;263;266;267 e 268, para as seguintes entidades:'
$ws.Rows.Item(521).RowHeight = 112

$ws.Range("A522").Value = 'svat_t35'
$ws.Range("B522").Value = 'Representação dos impostos diferidos no Ativo e no Passivo'
$ws.Rows.Item(522).RowHeight = 32

$ws.Range("A523").Value = 'svat_t35_ok'
$ws.Range("B523").Value = 'Verificamos que apenas tem saldo numa sub-conta da conta 274- "Impostos Diferidos". Se eventualmente efetuou compensação de saldos entre as contas 2741 e 2742, recomendamos verificar se essa compensação é adequada tendo em conta os parágrafos 68 a 69 da NCRF 25.'
$ws.Rows.Item(523).RowHeight = 96

$ws.Range("A524").Value = 'svat_t35_nok'
$ws.Range("B524").Value = 'Verificamos que tem saldos nas contas 2741 -''Ativos por impostos diferidos'' e 2742-''Passivos por impostos diferidos''. Sugerimos rever se não deverá efetuar compensação de saldos de acordo com os parágrafos 68 a 69 da NCRF 25.'
$ws.Rows.Item(524).RowHeight = 80

# Resize the i18n table to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G539"))

# Update worksheet view/selection to match the saved state
$excel.ActiveWindow.ScrollRow = 517
$ws.Range("A520").Select()

Write-Host "done"